# Updated cryptos list on Sat Dec  2 10:06:50 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# tracker sheet with the latest scraped values. Price strings are written
# through a text-formatted round trip (NumberFormat "@" -> Value -> Style
# "Normal") so Excel keeps digit-grouped price strings like "38.789.21"
# and trailing-zero values like "1.00" exactly as text instead of silently
# coercing them into numbers (which would also mangle the thousand-dot
# formatted prices). The cell keeps its original default/general style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-PriceText "D2" "38.789.21"
$ws.Range("E2").Value = "  +0.05%  "
# Row 3
Set-PriceText "D3" "2.103.49"
$ws.Range("E3").Value = "  -0.06%  "
# Row 5
Set-PriceText "D5" "228.54"
$ws.Range("E5").Value = "  -0.50%  "
# Row 6
$ws.Range("E6").Value = "  +0.07%  "
# Row 7
Set-PriceText "D7" "62.37"
$ws.Range("E7").Value = "  +0.79%  "
# Row 8
Set-PriceText "D8" "1.00"
$ws.Range("E8").Value = "  +0.04%  "
# Row 9
$ws.Range("E9").Value = "  +2.03%  "
# Row 10
$ws.Range("E10").Value = "  -0.61%  "
# Row 11
$ws.Range("E11").Value = "  -1.38%  "
# Row 12
Set-PriceText "D12" "15.76"
$ws.Range("E12").Value = "  +6.22%  "
# Row 13
Set-PriceText "D13" "2.415.40"
$ws.Range("E13").Value = "  -0.02%  "
# Row 14
Set-PriceText "D14" "22.12"
$ws.Range("E14").Value = "  -1.80%  "
# Row 15
$ws.Range("E15").Value = "  +3.43%  "
# Row 16
Set-PriceText "D16" "5.52"
$ws.Range("E16").Value = "  +0.48%  "
# Row 17
Set-PriceText "D17" "2.108.11"
$ws.Range("E17").Value = "  +0.18%  "
# Row 18
Set-PriceText "D18" "38.742.85"
$ws.Range("E18").Value = "  +0.26%  "
# Row 19
$ws.Range("E19").Value = "  +1.13%  "
# Row 20
Set-PriceText "D20" "6.11"
$ws.Range("E20").Value = "  +0.60%  "
# Row 21
$ws.Range("E21").Value = "  +0.49%  "
# Row 22
Set-PriceText "D22" "228.50"
$ws.Range("E22").Value = "  +0.66%  "
# Row 23
$ws.Range("E23").Value = "  +0.00%  "
# Row 24
Set-PriceText "D24" "2.34"
$ws.Range("E24").Value = "  -3.60%  "
# Row 25
$ws.Range("E25").Value = "  -0.09%  "
# Row 26
Set-PriceText "D26" "9.63"
$ws.Range("E26").Value = "  +1.72%  "
# Row 27
Set-PriceText "D27" "172.03"
# Row 28
Set-PriceText "D28" "0.139"
$ws.Range("E28").Value = "  +5.89%  "
# Row 29
$ws.Range("E29").Value = "  +4.35%  "
# Row 30
$ws.Range("E30").Value = "  +0.76%  "
# Row 31
Set-PriceText "D31" "2.53"
$ws.Range("E31").Value = "  +10.20%  "
# Row 32
$ws.Range("E32").Value = "  +0.34%  "
# Row 33
$ws.Range("E33").Value = "  +1.52%  "
# Row 34
$ws.Range("E34").Value = "  -0.82%  "
# Row 35
$ws.Range("E35").Value = "  +5.55%  "
# Row 36
Set-PriceText "D36" "0.0619"
$ws.Range("E36").Value = "  +1.85%  "
# Row 37
Set-PriceText "D37" "2.42"
$ws.Range("E37").Value = "  +0.40%  "
# Row 38
$ws.Range("E38").Value = "  +1.14%  "
# Row 39
$ws.Range("E39").Value = "  +0.08%  "
# Row 40
$ws.Range("E40").Value = "  -3.73%  "
# Row 41
Set-PriceText "D41" "102.95"
$ws.Range("E41").Value = "  +2.51%  "
# Row 42
$ws.Range("E42").Value = "  +3.11%  "
# Row 43
Set-PriceText "D43" "1.532.22"
$ws.Range("E43").Value = "  -0.97%  "
# Row 44
Set-PriceText "D44" "1.18"
$ws.Range("E44").Value = "  +5.59%  "
# Row 45
Set-PriceText "D45" "7.84"
$ws.Range("E45").Value = "  +2.50%  "
# Row 46
$ws.Range("E46").Value = "  -1.05%  "
# Row 47
$ws.Range("E47").Value = "  -0.58%  "
# Row 48
Set-PriceText "D48" "4.13"
$ws.Range("E48").Value = "  -1.08%  "
# Row 49
$ws.Range("E49").Value = "  +1.00%  "
# Row 50
$ws.Range("E50").Value = "  -0.31%  "
# Row 51
Set-PriceText "D51" "2.301.45"
$ws.Range("E51").Value = "  +0.05%  "
